$wb = $excel.ActiveWorkbook

# --- Sheet "Données brutes": rework AD (prime totale) & AF (bonus) columns ---
$ws = $wb.Worksheets.Item("Données brutes")
$ws.Activate()

for ($r = 9; $r -le 25; $r++) {
    $adCell = $ws.Range("AD$r")
    $adFormula = '=_xlfn.IFS(W' + $r + '="Oui",((J' + $r + '-O' + $r + ')+(I' + $r + '-N' + $r + ')*(P' + $r + '/100)),(AND(W' + $r + '="Oui",Y' + $r + '="Oui")),0)'
    $adCell.FormulaArray = $adFormula
    $adCell.NumberFormat = '"$"#,##0.00'
    $adCell.HorizontalAlignment = -4108

    $afCell = $ws.Range("AF$r")
    $afCell.Formula = '=IF(AE' + $r + '="Oui",(15000/8)/(N' + $r + '/O' + $r + '),0)'
    $afCell.NumberFormat = '"$"#,##0.00'
    $afCell.HorizontalAlignment = -4108
}

# --- View changes ---
$ws.Range("AI8").Select()
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1

$ws1 = $wb.Worksheets.Item("Résultats attendus")
$ws1.Activate()
$excel.ActiveWindow.Zoom = 96

$ws.Activate()
